$p = $ppt.ActivePresentation
$s = $p.Slides.Add(11, 1)
$notes = $s.NotesPage
$shp2 = $notes.Shapes.Item(2)
$tr = $shp2.TextFrame.TextRange
$tr.Text = "Hello"
Write-Host "After set: '$($tr.Text)' len=[$($tr.Length)]"
